# TX Directive Template: tidy up proofing marks (merge split runs that were
# only fragmented by stray grammar/spell-check markers) and rename the
# client city/county merge fields to the shorter city/county fields.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output ("NOT FOUND: " + $old)
    }
}

# Intro paragraph: drop the gramStart/gramEnd split around "In particular, you"
Replace-Text `
    "This is an important legal document known as an Advance Directive.  It is designed to help you communicate your wishes about medical treatment at some time in the future when you are unable to make your wishes known because of illness or injury.  These wishes are usually based on personal values.  In particular, you may want to consider what burdens or hardships of treatment you would be willing to accept for a particular amount of benefit obtained if you were seriously ill." `
    "This is an important legal document known as an Advance Directive.  It is designed to help you communicate your wishes about medical treatment at some time in the future when you are unable to make your wishes known because of illness or injury.  These wishes are usually based on personal values.  In particular, you may want to consider what burdens or hardships of treatment you would be willing to accept for a particular amount of benefit obtained if you were seriously ill."

# "I, {clientName}, recognize ... together as long as I am ... honored:" -
# remove gramStart/gramEnd split around "as / long as"
Replace-Text `
    ", recognize that the best health care is based upon a partnership of trust and communication with my physician.  My physician and I will make health care or treatment decisions together as long as I am of sound mind and able to make my wishes known.  If there comes a time that I am unable to make medical decisions about myself because of illness or injury, I direct that the following treatment preferences be honored:" `
    ", recognize that the best health care is based upon a partnership of trust and communication with my physician.  My physician and I will make health care or treatment decisions together as long as I am of sound mind and able to make my wishes known.  If there comes a time that I am unable to make medical decisions about myself because of illness or injury, I direct that the following treatment preferences be honored:"

# "...die as gently as possible; " appears twice - replaces both occurrences
Replace-Text `
    "I request that all treatments other than those needed to keep me comfortable be discontinued or withheld and my physician allow me to die as gently as possible; " `
    "I request that all treatments other than those needed to keep me comfortable be discontinued or withheld and my physician allow me to die as gently as possible; "

# Additional requests explanatory text
Replace-Text `
    "(After discussion with your physician, you may wish to consider listing particular treatments in this space that you do or do not want in specific circumstances, such as artificially administered nutrition and hydration, intravenous antibiotics, etc.  Be sure to state whether you do or do not want the particular treatment.)" `
    "(After discussion with your physician, you may wish to consider listing particular treatments in this space that you do or do not want in specific circumstances, such as artificially administered nutrition and hydration, intravenous antibiotics, etc.  Be sure to state whether you do or do not want the particular treatment.)"

# City/County merge fields
Replace-Text "clientCity" "city"
Replace-Text "clientCounty" "county"

# Notary county merge-field run cleanup
Replace-Text `
    "{#notaryCounty}{notaryCounty}{/notaryCounty}{^notaryCount" `
    "{#notaryCounty}{notaryCounty}{/notaryCounty}{^notaryCount"

# "Irreversible condition" bullet point
Replace-Text `
    "that may be treated, but is never cured or eliminated;" `
    "that may be treated, but is never cured or eliminated;"

# Explanation paragraph (irreversible)
Replace-Text `
    "Explanation:  Many serious illnesses such as cancer, failure of major organs (kidney, heart, liver, or lung), and serious brain disease such as Alzheimer’s dementia may be considered irreversible early on.  There is no cure, but the patient may be kept alive for prolonged periods of time if the patient receives life-sustaining treatments.  Late in the course of the same illness, the disease may be considered terminal when, even with treatment, the patient is expected to die.  You may wish to consider which burdens of treatment you would be willing to accept in an effort to achieve a particular outcome.  This is a very personal decision that you may wish to discuss with your physician, family, or other important persons in your life." `
    "Explanation:  Many serious illnesses such as cancer, failure of major organs (kidney, heart, liver, or lung), and serious brain disease such as Alzheimer’s dementia may be considered irreversible early on.  There is no cure, but the patient may be kept alive for prolonged periods of time if the patient receives life-sustaining treatments.  Late in the course of the same illness, the disease may be considered terminal when, even with treatment, the patient is expected to die.  You may wish to consider which burdens of treatment you would be willing to accept in an effort to achieve a particular outcome.  This is a very personal decision that you may wish to discuss with your physician, family, or other important persons in your life."

# Explanation paragraph (terminal)
Replace-Text `
    "Explanation:  Many serious illnesses may be considered irreversible early in the course of the illness, but they may not be considered terminal until the disease is fairly advanced.  In thinking about terminal illness and its treatment, you again may wish to consider the relative benefits and burdens of treatment and discuss your wishes with your physician, family, or other important persons in your life." `
    "Explanation:  Many serious illnesses may be considered irreversible early in the course of the illness, but they may not be considered terminal until the disease is fairly advanced.  In thinking about terminal illness and its treatment, you again may wish to consider the relative benefits and burdens of treatment and discuss your wishes with your physician, family, or other important persons in your life."

Write-Output "done"
